$d = $word.ActiveDocument

# 1. Task 1 paragraph: merge "In task 1 the APK files for the " + "Bitmoji" + " Android app was downloaded..."
$d.Content.Find.Execute(
    "In task 1 the APK files for the Bitmoji Android app was downloaded. This task showed how easy it was to obtain APK files for Android apps to perform the repackaging attack. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In task 1 the APK files for the Bitmoji Android app was downloaded. This task showed how easy it was to obtain APK files for Android apps to perform the repackaging attack. ",
    2)

# 2. Task 2 first part: merge "In task 2 we used " + "APKTool" + " to "
$d.Content.Find.Execute(
    "In task 2 we used APKTool to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In task 2 we used APKTool to ",
    2)

# 3. Task 2 second part: merge " the " + "Bitmoji" + " " + "dex" + " code to " + "smali" + " code. ... Android apps. "
$d.Content.Find.Execute(
    " the Bitmoji dex code to smali code. In this task it is easy to see that if you have APKTool it is easy to disassemble any APK files for Android apps. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " the Bitmoji dex code to smali code. In this task it is easy to see that if you have APKTool it is easy to disassemble any APK files for Android apps. ",
    2)

# 4. Task 3 paragraph: merge "In task 3 we inserted a " + "MaliciousCode.smali" + " file in the " + "smali" + "/com folder ..."
$d.Content.Find.Execute(
    "In task 3 we inserted a MaliciousCode.smali file in the smali/com folder and then edited the AndroidManifest.xml file so that ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In task 3 we inserted a MaliciousCode.smali file in the smali/com folder and then edited the AndroidManifest.xml file so that ",
    2)

# 5. Task 4 paragraph part A: merge "...repackage the application with the " + "APKTool" + ". Once the "
$d.Content.Find.Execute(
    "Since we have the malicious code is inserted now we have to repackage the application with the APKTool. Once the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Since we have the malicious code is inserted now we have to repackage the application with the APKTool. Once the ",
    2)

# 5b. Task 4 paragraph part B: merge "...using the " + "keytool" + " command. Then " + "jarsigner" + " is used to sign..."
$d.Content.Find.Execute(
    "application is repackaged we have to sign the APK file so generated the public and private key using the keytool command. Then jarsigner is used to sign the APK file using the keys that were generated. In this task we see that even though we are the ones inserting the malicious code we still need to make sure no one else can get to the code we inserted or the purpose of inserting it would be pointless.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "application is repackaged we have to sign the APK file so generated the public and private key using the keytool command. Then jarsigner is used to sign the APK file using the keys that were generated. In this task we see that even though we are the ones inserting the malicious code we still need to make sure no one else can get to the code we inserted or the purpose of inserting it would be pointless.",
    2)

# 6. Task 5 paragraph: merge "The last part..." + "Bitmoji" + " app with the malicious code..." and
#    change "...deleting the contacts in " + "Bitmoji" + " app. " to "...deleting the contacts in the VM. "
$d.Content.Find.Execute(
    "The last part of the lab has us install the Bitmoji app with the malicious code onto the Android VM and see if the attack works. With this task we see that our malicious code is successful in deleting the contacts in Bitmoji app. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The last part of the lab has us install the Bitmoji app with the malicious code onto the Android VM and see if the attack works. With this task we see that our malicious code is successful in deleting the contacts in the VM. ",
    2)

# 7. Question 3 answer paragraph: merge "...antivirus " + "softwares" + " are not always 100% effective. ..."
$d.Content.Find.Execute(
    "No, Google Play Store cannot totally keep a customer safe from the attacks, much in the same way that antivirus softwares are not always 100% effective. Even if Google was to deploy a built-in malware scanner, it would have to be updated constantly in order to keep up with the newest malware.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "No, Google Play Store cannot totally keep a customer safe from the attacks, much in the same way that antivirus softwares are not always 100% effective. Even if Google was to deploy a built-in malware scanner, it would have to be updated constantly in order to keep up with the newest malware.",
    2)
